$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape "CaixaDeTexto 8" (INTRODUÇÃO text box)
# ---------------------------------------------------------------------
$sh9 = $s.Shapes.Item("CaixaDeTexto 8")

# Resize the textbox height (cy 2862322 -> 2492990 EMU, i.e. 225.3796850... -> 196.2984251968504 pt)
$sh9.Height = 196.2984251968504

$tr9 = $sh9.TextFrame.TextRange

# Second paragraph's run text (offset 13, length 149 chars in the original text)
$c9p2 = $tr9.Characters(13, 149)
$c9p2.Text = "`tA utilidade da localização de pessoas e aparelhos vem se aperfeiçoando cada vez mais com a evolução de programas e aplicativos que utilizam-se deste recurso, seja para trazer mais mobilidade no trânsito ou para simplificar a procura de um restaurante mais próximo, por exemplo."

# Third paragraph's run text. Re-resolve the paragraph start after the edit above
# (its own length changed), so the offset used here stays valid.
$tr9b = $sh9.TextFrame.TextRange
$para3 = $tr9b.Paragraphs(3, 1)
$c9p3 = $tr9b.Characters($para3.Start, 593)
$c9p3.Text = "`tCom a baixa efetividade da localização por GPS em lugares fechados, somada a evolução dos microcontroladores, a localização indoor acaba se tornando alvo de pesquisas e projetos por ser bastante relevante em ambientes grandes e fechados onde a vida de muitas pessoas seria facilitada por ter como se situar, e onde seria possível monitorar funcionários dentro de uma empresa. "

# ---------------------------------------------------------------------
# Shape "CaixaDeTexto 11" (RESULTADOS E DISCUSSÕES text box)
# ---------------------------------------------------------------------
$sh8 = $s.Shapes.Item("CaixaDeTexto 11")
$tr8 = $sh8.TextFrame.TextRange
$c8 = $tr8.Characters(26, 553)
$c8.Text = "A partir dos desafios propostos, apenas sabemos a localização da sala e o andar que o funcionário se encontra, embora alguns ajustes ainda necessitam ser realizados. Os dados captados pelo protótipo, pode ser consultado por um dispositivo qualquer com acesso a internet e é possível saber onde a pessoa se encontra dentro da sala, embora apresentando um grande erro relacionado a localização do mesmo. Portanto, necessitaria da análise da onda emitida pela esp8266 e, junto à isto, verificar como efetuar leituras precisas a fim de gerar uma localização mais exata."

# ---------------------------------------------------------------------
# Shape "CaixaDeTexto 13" (CONSIDERAÇÕES FINAIS text box)
# Edit the later-offset run first so the earlier run's character offset
# (which precedes it in the text stream) is not invalidated by the
# text-length change of the later edit.
# ---------------------------------------------------------------------
$sh10 = $s.Shapes.Item("CaixaDeTexto 13")
$tr10 = $sh10.TextFrame.TextRange
$c10b = $tr10.Characters(194, 220)
$c10b.Text = " devido a fatores externos e análise do sinal gerado pela esp8266. Contudo, o protótipo é capaz de localizar a sala e o andar que a pessoa estaria e com isto pode ser identificado se a pessoa é ou não funcionário da empresa."

$tr10b = $sh10.TextFrame.TextRange
$c10a = $tr10b.Characters(23, 164)
$c10a.Text = "A partir dos dados gerados pelo protótipo é possível analisar que a precisão da localização de uma pessoa fica fora do que seria considerado razoável para a localização "
